# Weekly fruit/vegetable (Betarraga) data update:
# Two new daily records (for a new date, serial 44606) are inserted at row 150,
# pushing the existing rows 150-248 down by two rows (to 152-250).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at position 150, shifting rows 150:248 down to 152:250.
$ws.Range("A150:R151").Insert()

# New row 150 ("Primera" quality)
$ws.Cells.Item(150,1).Value = 1
$ws.Cells.Item(150,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(150,3).Value = "Arica y Parinacota"
$ws.Cells.Item(150,4).Value = 44606
$ws.Cells.Item(150,5).Value = 15
$ws.Cells.Item(150,6).Value = 100114014
$ws.Cells.Item(150,7).Value = "Betarraga"
$ws.Cells.Item(150,8).Value = "Sin especificar"
$ws.Cells.Item(150,9).Value = "Primera"
$ws.Cells.Item(150,10).Value = 1000
$ws.Cells.Item(150,11).Value = 400
$ws.Cells.Item(150,12).Value = 450
$ws.Cells.Item(150,13).Value = 425
$ws.Cells.Item(150,14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(150,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(150,16).Value = 106
$ws.Cells.Item(150,17).Value = 4
$ws.Cells.Item(150,18).Value = "Hortaliza"

# New row 151 ("Segunda" quality)
$ws.Cells.Item(151,1).Value = 1
$ws.Cells.Item(151,2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(151,3).Value = "Arica y Parinacota"
$ws.Cells.Item(151,4).Value = 44606
$ws.Cells.Item(151,5).Value = 15
$ws.Cells.Item(151,6).Value = 100114014
$ws.Cells.Item(151,7).Value = "Betarraga"
$ws.Cells.Item(151,8).Value = "Sin especificar"
$ws.Cells.Item(151,9).Value = "Segunda"
$ws.Cells.Item(151,10).Value = 1000
$ws.Cells.Item(151,11).Value = 400
$ws.Cells.Item(151,12).Value = 450
$ws.Cells.Item(151,13).Value = 425
$ws.Cells.Item(151,14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(151,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(151,16).Value = 85
$ws.Cells.Item(151,17).Value = 5
$ws.Cells.Item(151,18).Value = "Hortaliza"
